$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Paragraph 3 ("Shall type be global variables ?") is retyped as
#    four separate runs: "T" / "ype" / " variables are" / " global
#    variables" (-> "Type variables are global variables").
#    The paragraph keeps its ListParagraph/numPr formatting.
#    NOTE: the existing "_GoBack" bookmark is left in place while we
#    do this (deleting it first causes the freshly-inserted runs to
#    coalesce into one run); it gets relocated afterwards.
# -----------------------------------------------------------------
$p3 = $d.Paragraphs(3)

$start = $p3.Range.Start
$clearRange = $d.Range($start, $p3.Range.End - 1)
$clearRange.Text = ""

$pos = $start
$chunks = @("T", "ype", " variables are", " global variables")
foreach ($chunk in $chunks) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($chunk)
    $pos = $pos + $chunk.Length
}

# -----------------------------------------------------------------
# 2) Insert a brand-new plain paragraph right after it, holding the
#    "15/10: ..." note plus the underlined "loop remains".
# -----------------------------------------------------------------
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$p4.Range.ListFormat.RemoveNumbers()
$p4.Style = "Normal"

# Relocate the "_GoBack" bookmark out of paragraph 3 now that its
# runs are final.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$pos2 = $p4.Range.Start
$r1 = $d.Range($pos2, $pos2)
$r1.InsertAfter("15/10: Parts 1 and 3 complete, ")
$pos2 = $pos2 + "15/10: Parts 1 and 3 complete, ".Length

$bmStart = $pos2
$r2 = $d.Range($pos2, $pos2)
$r2.InsertAfter("loop remains")
$r2.Font.Underline = 1
$pos2 = $pos2 + "loop remains".Length

# Bookmark wraps the "loop remains" run, same as it used to wrap the
# (now removed) trailing text of paragraph 3.
$bmRange = $d.Range($bmStart, $pos2)
$d.Bookmarks.Add("_GoBack", $bmRange)
